# "Generate Report for Handback" — mark the two e2e files as handed back,
# fill in the per-language "Latest Target File" / "Latest Handback File"
# columns with links + filenames, stamp a handback datetime, and widen the
# columns that now hold longer file names.

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45e569324548f8ac74f71497c93f0f4b8db501fe/e2e/"
$file1 = "3ef57ec2-dcc2-4626-b272-9dcfe06555c4.md"
$file2 = "6d3486c6-7daa-4f36-a7c4-2158e56e2baf.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both status cells per row (zh-cn + de-de columns)
# flip from "Ready for handoff" to "Handed back: in sync with en-US".
# Also widen columns E/F to fit the longer text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("C3").Value = $statusHandedBack

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlBase + $file1, "", "", $file1)
$wsZh.Range("J2").Value = "3ef57ec2-dcc2-4626-b272-9dcfe06555c4.fa296af86d6eca9396d0803dd4bfa5a5687fc8af.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlBase + $file2, "", "", $file2)
$wsZh.Range("J3").Value = "6d3486c6-7daa-4f36-a7c4-2158e56e2baf.0b615aecdb7ab1922a17e1984941a43019087f44.zh-cn.xlf"

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("C3").Value = $statusHandedBack

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlBase + $file1, "", "", $file1)
$wsDe.Range("J2").Value = "3ef57ec2-dcc2-4626-b272-9dcfe06555c4.fa296af86d6eca9396d0803dd4bfa5a5687fc8af.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 02:31:58"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlBase + $file2, "", "", $file2)
$wsDe.Range("J3").Value = "6d3486c6-7daa-4f36-a7c4-2158e56e2baf.0b615aecdb7ab1922a17e1984941a43019087f44.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-30 02:31:58"

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
